$d = $word.ActiveDocument

# Locate the target paragraph precisely by its distinctive original text,
# rather than relying on a hard-coded paragraph index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "Citizens must not vote in any other Wahlbezirk than the one they are registered in xor by Briefwahl.*") {
        $target = $p
        break
    }
}

$pStart = $target.Range.Start

# Offsets below are relative to the paragraph start, computed from the
# known original run boundaries:
#   "Citizens "                               ->  0 .. 9
#   "must"                                     ->  9 .. 13
#   " not vote in any other "                  -> 13 .. 36
#   "Wahlbezirk"                               -> 36 .. 46
#   " than the one they are registered in "    -> 46 .. 83
#   "x"                                        -> 83 .. 84
#   "or"                                       -> 84 .. 86
#   " by "                                     -> 86 .. 90
#   "Briefwahl"                                -> 90 .. 99
#   "."                                        -> 99 .. 100

# --- Step 1 (rightmost edit first, so earlier offsets stay valid) ---
# " than the one they are registered in " -> " they are registered " + in ""
# with the _GoBack bookmark relocated to the split point between them.
$r5 = $d.Range($pStart + 46, $pStart + 83)
$r5.Text = " they are registered in "

$bmPos = $pStart + 46 + (" they are registered ").Length
$bmRng = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

# --- Step 2 ---
# " not vote in any other " -> " only be possible in the one "
$r3 = $d.Range($pStart + 13, $pStart + 36)
$r3.Text = " only be possible in the one "

# --- Step 3 (leftmost edit last) ---
# "Citizens " -> "Voting of citizens "
$r1 = $d.Range($pStart + 0, $pStart + 9)
$r1.Text = "Voting of citizens "

Write-Output ("Final paragraph text: [" + $target.Range.Text + "]")
